$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    # Force text so a numeric-looking string (e.g. "337.10") is stored
    # verbatim as text rather than auto-coerced into a Number (which would
    # silently drop meaningful trailing/format digits), matching the
    # original inline-string cell. Reset the style back to Normal afterwards
    # so we do not leave a stray text-format style on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '28.107.01'
Set-TextValue 'E2' '  +0.17%  '

# Row 3
Set-TextValue 'D3' '1.797.04'
Set-TextValue 'E3' '  +2.37%  '

# Row 4
Set-TextValue 'D4' '1.006'
Set-TextValue 'E4' '  +0.17%  '

# Row 5
Set-TextValue 'D5' '337.10'
Set-TextValue 'E5' '  -0.03%  '

# Row 6
Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  +0.26%  '

# Row 7
Set-TextValue 'D7' '0.4600'
Set-TextValue 'E7' '  +21.81%  '

# Row 8
Set-TextValue 'D8' '0.3692'
Set-TextValue 'E8' '  +10.15%  '

# Row 9
Set-TextValue 'D9' '45.27'
Set-TextValue 'E9' '  -0.39%  '

# Row 10
Set-TextValue 'D10' '0.07645'
Set-TextValue 'E10' '  +6.19%  '

# Row 11
Set-TextValue 'D11' '1.144'
Set-TextValue 'E11' '  +2.37%  '

# Row 12
Set-TextValue 'B12' 'BinanceUSD'
Set-TextValue 'C12' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D12' '1.004'
Set-TextValue 'E12' '  +0.21%  '

# Row 13
Set-TextValue 'B13' 'Solana'
Set-TextValue 'C13' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D13' '22.50'
Set-TextValue 'E13' '  -0.20%  '

# Row 14
Set-TextValue 'D14' '6.335'
Set-TextValue 'E14' '  +3.16%  '

# Row 15
Set-TextValue 'D15' '7.384'
Set-TextValue 'E15' '  +3.42%  '

# Row 16
Set-TextValue 'D16' '1.794.70'
Set-TextValue 'E16' '  +2.10%  '

# Row 17
Set-TextValue 'D17' '0.00001090'
Set-TextValue 'E17' '  +3.17%  '

# Row 18
Set-TextValue 'D18' '0.06717'
Set-TextValue 'E18' '  +1.99%  '

# Row 19
Set-TextValue 'D19' '82.37'
Set-TextValue 'E19' '  +2.30%  '

# Row 20
Set-TextValue 'D20' '1.002'
Set-TextValue 'E20' '  +0.13%  '

# Row 21
Set-TextValue 'D21' '17.39'
Set-TextValue 'E21' '  +3.04%  '

# Row 22
Set-TextValue 'D22' '6.399'
Set-TextValue 'E22' '  +2.58%  '

# Row 23
Set-TextValue 'D23' '28.108.68'
Set-TextValue 'E23' '  +0.09%  '

# Row 24
Set-TextValue 'D24' '11.87'
Set-TextValue 'E24' '  +1.83%  '

# Row 25
Set-TextValue 'D25' '2.413'
Set-TextValue 'E25' '  +0.73%  '

# Row 26
Set-TextValue 'D26' '20.73'
Set-TextValue 'E26' '  +4.69%  '

# Row 27
Set-TextValue 'D27' '2.384'
Set-TextValue 'E27' '  +2.77%  '

# Row 28
Set-TextValue 'D28' '151.55'
Set-TextValue 'E28' '  -0.72%  '

# Row 29
Set-TextValue 'D29' '2.003.17'
Set-TextValue 'E29' '  +2.40%  '

# Row 30
Set-TextValue 'D30' '133.62'
Set-TextValue 'E30' '  +1.41%  '

# Row 31
Set-TextValue 'D31' '1.256'
Set-TextValue 'E31' '  +0.72%  '

# Row 32
Set-TextValue 'D32' '4.055'
Set-TextValue 'E32' '  +0.87%  '

# Row 33
Set-TextValue 'D33' '0.09642'
Set-TextValue 'E33' '  +10.62%  '

# Row 34
Set-TextValue 'D34' '5.905'
Set-TextValue 'E34' '  +2.09%  '

# Row 35
Set-TextValue 'D35' '0.02373'
Set-TextValue 'E35' '  +2.06%  '

# Row 36
Set-TextValue 'D36' '0.2219'
Set-TextValue 'E36' '  +5.15%  '

# Row 37
Set-TextValue 'D37' '12.15'
Set-TextValue 'E37' '  -0.45%  '

# Row 38
Set-TextValue 'B38' 'TheSandbox'
Set-TextValue 'C38' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D38' '0.6699'
Set-TextValue 'E38' '  +0.24%  '

# Row 39
Set-TextValue 'B39' 'Hedera'
Set-TextValue 'C39' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D39' '0.06328'
Set-TextValue 'E39' '  +2.32%  '

# Row 40
Set-TextValue 'D40' '5.252'
Set-TextValue 'E40' '  +1.86%  '

# Row 41
Set-TextValue 'D41' '1.505'
Set-TextValue 'E41' '  +4.06%  '

# Row 42
Set-TextValue 'D42' '1.235'
Set-TextValue 'E42' '  +1.53%  '

# Row 43
Set-TextValue 'D43' '8.066'
Set-TextValue 'E43' '  +0.55%  '

# Row 44
Set-TextValue 'E44' '  +3.95%  '

# Row 45
Set-TextValue 'E45' '  +0.15%  '

# Row 46
Set-TextValue 'D46' '0.6147'
Set-TextValue 'E46' '  +1.61%  '

# Row 47
Set-TextValue 'D47' '3.846'
Set-TextValue 'E47' '  +0.29%  '

# Row 48
Set-TextValue 'D48' '130.10'
Set-TextValue 'E48' '  +1.07%  '

# Row 49
Set-TextValue 'D49' '2.052'
Set-TextValue 'E49' '  +1.76%  '

# Row 50
Set-TextValue 'D50' '1.181'
Set-TextValue 'E50' '  +0.99%  '

# Row 51
Set-TextValue 'D51' '0.07131'
Set-TextValue 'E51' '  -0.29%  '
